# Apply the commit's edit to the "硫酸(折100%)" (sulfuric acid) worksheet:
#  1. For each year block of 4 rows (A/B/C/D sub-periods), the "B" row and
#     "C" row have their data swapped (the B-labelled row's figures move to
#     the row that used to hold the C-labelled row's figures, and vice
#     versa) while columns A-E keep their content attached to the row it
#     belongs with.
#  2. Columns F ("...产销率") and G ("...销售量") are removed entirely
#     (including their header cells in row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs (1-based sheet rows) whose A:E contents must be swapped.
# Pattern: starting at row 3, every 4-row year block swaps its 2nd and 3rd
# rows (rows 3&4, 7&8, 11&12, ... 67&68).
for ($k = 0; $k -lt 17; $k++) {
    $row1 = 4 * $k + 3
    $row2 = 4 * $k + 4

    $rng1 = $ws.Range("A" + $row1 + ":E" + $row1)
    $rng2 = $ws.Range("A" + $row2 + ":E" + $row2)

    $vals1 = $rng1.Value2
    $vals2 = $rng2.Value2

    $rng1.Value2 = $vals2
    $rng2.Value2 = $vals1
}

# Remove columns F:G (and their header cells) entirely; Excel shifts the
# dimension / used range down to A1:E69 automatically.
$ws.Range("F1:G1").EntireColumn.Delete()
